$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "(203957296, Omri Ben Shabat: -2,-9)"
$ws.Range("B1").Value = "(206532695, Matan Vakrat: 8,-6)"
$ws.Range("C1").Value = "(302962915, Asher  Odeh: -7,-10)"
$ws.Range("D1").Value = "(308035542, Anastasia  Kubi: -8,3)"
$ws.Range("E1").Value = "(311177802, Christina  Uksusman: -8,5)"
$ws.Range("F1").Value = "(305251175, Or  Leder: -7,-7)"

$ws.Range("A3").Value = "cost: 561.8412363853096"
$ws.Range("A4").Value = "time: 67.1051545481637"
